# Auto-generated edit script to apply the diff to before.xlsx
$wb = $excel.ActiveWorkbook

# ======== Sheet: ALC ========
$ws = $wb.Worksheets.Item("ALC")
# @@ -715,25 +715,25 @@
$ws.Range("H2").Value = 1531.2142
$ws.Range("I2").Value = 651.2857
$ws.Range("J2").Value = 2411.1428
$ws.Range("K2").Value = 651.2857
$ws.Range("L2").Value = 2411.1428
$ws.Range("M2").Value = -538.2857
$ws.Range("N2").Value = -2637.1428
# @@ -4653,22 +4653,22 @@
$ws.Range("H82").Value = 750
$ws.Range("I82").Value = 750
$ws.Range("K82").Value = 2250
$ws.Range("M82").Value = -1844
# @@ -4800,22 +4800,22 @@
$ws.Range("H85").Value = 750
$ws.Range("I85").Value = 750
$ws.Range("K85").Value = 2250
$ws.Range("M85").Value = -846
# @@ -6337,22 +6337,22 @@
$ws.Range("H116").Value = 3848.8572
$ws.Range("I116").Value = 3496
$ws.Range("K116").Value = 3496
$ws.Range("M116").Value = -54
# @@ -7103,25 +7103,22 @@
$ws.Range("H132").Value = 1746.2222
$ws.Range("I132").Value = 1746.2222
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5238.6666
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2708.6666
$ws.Range("N132").Value = $null
# @@ -7394,25 +7391,25 @@
$ws.Range("H138").Value = 3219.6758
$ws.Range("I138").Value = 1632
$ws.Range("J138").Value = 4079.6667
$ws.Range("K138").Value = 4896
$ws.Range("L138").Value = 12239.0001
$ws.Range("M138").Value = 244
$ws.Range("N138").Value = -22519.0001

# ======== Sheet: ARM ========
$ws = $wb.Worksheets.Item("ARM")
# @@ -9797,25 +9794,25 @@
$ws.Range("H45").Value = 5288.8
$ws.Range("I45").Value = 5481.5
$ws.Range("J45").Value = 4999.75
$ws.Range("K45").Value = 5481.5
$ws.Range("L45").Value = 4999.75
$ws.Range("M45").Value = -5104.5
$ws.Range("N45").Value = -5753.75
# @@ -10143,19 +10140,22 @@
$ws.Range("H52").Value = 10709
$ws.Range("I52").Value = 10709
$ws.Range("K52").Value = 10709
$ws.Range("M52").Value = -10391
# @@ -10566,22 +10566,22 @@
$ws.Range("H61").Value = 7495
$ws.Range("I61").Value = 1658.3334
$ws.Range("K61").Value = 1658.3334
$ws.Range("M61").Value = -1446.3334
# @@ -11200,25 +11200,25 @@
$ws.Range("H74").Value = 909.875
$ws.Range("J74").Value = 876
$ws.Range("L74").Value = 876
$ws.Range("N74").Value = -2624
# @@ -11344,25 +11344,25 @@
$ws.Range("H77").Value = 909.875
$ws.Range("J77").Value = 876
$ws.Range("L77").Value = 4380
$ws.Range("N77").Value = -13116
# @@ -12490,25 +12490,25 @@
$ws.Range("H101").Value = 19999.5
$ws.Range("J101").Value = 19999
$ws.Range("L101").Value = 19999
$ws.Range("N101").Value = -26489
# @@ -13489,25 +13489,25 @@
$ws.Range("H122").Value = 23206.867
$ws.Range("I122").Value = 26161.77
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 78485.31
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -76035.31
$ws.Range("N122").Value = -16900
# @@ -14160,22 +14160,22 @@
$ws.Range("H136").Value = 7495
$ws.Range("I136").Value = 1658.3334
$ws.Range("K136").Value = 4975.0002
$ws.Range("M136").Value = -2425.0002

# ======== Sheet: BSM ========
$ws = $wb.Worksheets.Item("BSM")
# @@ -15225,25 +15225,22 @@
$ws.Range("H16").Value = 25000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 25000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 25000
$ws.Range("M16").Value = $null
$ws.Range("N16").Value = -25340
# @@ -19499,22 +19496,22 @@
$ws.Range("H105").Value = 2312.6
$ws.Range("I105").Value = 2312.6
$ws.Range("K105").Value = 2312.6
$ws.Range("M105").Value = -565.5999999999999
# @@ -19597,25 +19594,25 @@
$ws.Range("H107").Value = 3052.6667
$ws.Range("I107").Value = 3168.3333
$ws.Range("J107").Value = 2590
$ws.Range("K107").Value = 3168.3333
$ws.Range("L107").Value = 2590
$ws.Range("M107").Value = -1248.3333
$ws.Range("N107").Value = -6430
# @@ -20881,22 +20878,22 @@
$ws.Range("H134").Value = 7765.5
$ws.Range("I134").Value = 7544.8
$ws.Range("K134").Value = 22634.4
$ws.Range("M134").Value = -20099.4

# ======== Sheet: CRP ========
$ws = $wb.Worksheets.Item("CRP")
# @@ -22065,25 +22062,25 @@
$ws.Range("H16").Value = 5753.875
$ws.Range("J16").Value = 6999.75
$ws.Range("L16").Value = 6999.75
$ws.Range("N16").Value = -7573.75
# @@ -22788,25 +22785,25 @@
$ws.Range("H31").Value = 1124.25
$ws.Range("I31").Value = 756
$ws.Range("J31").Value = 1247
$ws.Range("K31").Value = 756
$ws.Range("L31").Value = 1247
$ws.Range("M31").Value = -461
$ws.Range("N31").Value = -1837
# @@ -22944,25 +22941,25 @@
$ws.Range("H34").Value = 1124.25
$ws.Range("I34").Value = 756
$ws.Range("J34").Value = 1247
$ws.Range("K34").Value = 756
$ws.Range("L34").Value = 1247
$ws.Range("M34").Value = -554
$ws.Range("N34").Value = -1651
# @@ -26087,7 +26084,7 @@
$ws.Range("H99").Value = 1159.6
# @@ -26479,25 +26476,19 @@
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = $null
$ws.Range("N107").Value = $null
# @@ -26770,25 +26761,25 @@
$ws.Range("H113").Value = 5753.875
$ws.Range("J113").Value = 6999.75
$ws.Range("L113").Value = 6999.75
$ws.Range("N113").Value = -11339.75
# @@ -27395,7 +27386,7 @@
$ws.Range("H126").Value = 1159.6

# ======== Sheet: CUL ========
$ws = $wb.Worksheets.Item("CUL")
# @@ -28518,22 +28509,22 @@
$ws.Range("H7").Value = 927.9286
$ws.Range("I7").Value = 915.9167
$ws.Range("K7").Value = 2747.7501
$ws.Range("M7").Value = -2635.7501
# @@ -28720,22 +28711,22 @@
$ws.Range("H11").Value = 1476.6428
$ws.Range("I11").Value = 2667.5715
$ws.Range("K11").Value = 8002.7145
$ws.Range("M11").Value = -7862.7145
# @@ -29467,25 +29458,25 @@
$ws.Range("H26").Value = 731.1667
$ws.Range("I26").Value = 829.6667
$ws.Range("J26").Value = 632.6667
$ws.Range("K26").Value = 2489.0001
$ws.Range("L26").Value = 1898.0001
$ws.Range("M26").Value = -2201.0001
$ws.Range("N26").Value = -2474.0001
# @@ -29617,25 +29608,22 @@
$ws.Range("H29").Value = 576
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 576
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 1728
$ws.Range("M29").Value = $null
$ws.Range("N29").Value = -2282
# @@ -31953,25 +31941,25 @@
$ws.Range("H76").Value = 16626.625
$ws.Range("J76").Value = 19333.334
$ws.Range("L76").Value = 58000.00199999999
$ws.Range("N76").Value = -58766.00199999999
# @@ -32106,25 +32094,25 @@
$ws.Range("H79").Value = 16626.625
$ws.Range("J79").Value = 19333.334
$ws.Range("L79").Value = 58000.00199999999
$ws.Range("N79").Value = -60652.00199999999
# @@ -33931,22 +33919,22 @@
$ws.Range("H116").Value = 2250
$ws.Range("I116").Value = 2250
$ws.Range("K116").Value = 6750
$ws.Range("M116").Value = -3308
# @@ -34718,22 +34706,22 @@
$ws.Range("H132").Value = 1959.2
$ws.Range("I132").Value = 1899.5
$ws.Range("K132").Value = 17095.5
$ws.Range("M132").Value = -14565.5

# ======== Sheet: GSM ========
$ws = $wb.Worksheets.Item("GSM")
# @@ -40166,25 +40154,25 @@
$ws.Range("H102").Value = 2184.9167
$ws.Range("I102").Value = 1831.9
$ws.Range("J102").Value = 3950
$ws.Range("K102").Value = 1831.9
$ws.Range("L102").Value = 3950
$ws.Range("M102").Value = -209.9000000000001
$ws.Range("N102").Value = -7194
# @@ -40687,22 +40675,22 @@
$ws.Range("H113").Value = 6999
$ws.Range("I113").Value = 1498
$ws.Range("K113").Value = 1498
$ws.Range("M113").Value = 672

# ======== Sheet: LTW ========
$ws = $wb.Worksheets.Item("LTW")
# @@ -43140,25 +43128,25 @@
$ws.Range("H22").Value = 812.7778
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 903
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 903
$ws.Range("M22").Value = -405
$ws.Range("N22").Value = -1493
# @@ -43376,25 +43364,25 @@
$ws.Range("H27").Value = 812.7778
$ws.Range("I27").Value = 700
$ws.Range("J27").Value = 903
$ws.Range("K27").Value = 700
$ws.Range("L27").Value = 903
$ws.Range("M27").Value = -593
$ws.Range("N27").Value = -1117
# @@ -44016,22 +44004,22 @@
$ws.Range("H40").Value = 2865.5833
$ws.Range("I40").Value = 2132.6
$ws.Range("K40").Value = 2132.6
$ws.Range("M40").Value = -1996.6
# @@ -44739,25 +44727,25 @@
$ws.Range("H55").Value = 1296.8
$ws.Range("I55").Value = 1133.75
$ws.Range("J55").Value = 1483.1428
$ws.Range("K55").Value = 1133.75
$ws.Range("L55").Value = 1483.1428
$ws.Range("M55").Value = -960.75
$ws.Range("N55").Value = -1829.1428
# @@ -45033,22 +45021,19 @@
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = $null
# @@ -47190,22 +47175,22 @@
$ws.Range("H106").Value = 20497
$ws.Range("J106").Value = 20497
$ws.Range("L106").Value = 20497
$ws.Range("N106").Value = -23021
# @@ -47524,22 +47509,19 @@
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = $null
# @@ -48437,25 +48419,25 @@
$ws.Range("H132").Value = 4659.4
$ws.Range("I132").Value = 3699.5
$ws.Range("J132").Value = 5299.3335
$ws.Range("K132").Value = 11098.5
$ws.Range("L132").Value = 15898.0005
$ws.Range("M132").Value = -8568.5
$ws.Range("N132").Value = -20958.0005

# ======== Sheet: WVR ========
$ws = $wb.Worksheets.Item("WVR")
# @@ -49122,22 +49104,22 @@
$ws.Range("H4").Value = 7860999.5
$ws.Range("I4").Value = 16670666
$ws.Range("K4").Value = 16670666
$ws.Range("M4").Value = -16670553
# @@ -52556,23 +52538,26 @@
$ws.Range("H75").Value = 90039.336
$ws.Range("J75").Value = 90000
$ws.Range("L75").Value = 90000
$ws.Range("N75").Value = -91872
# @@ -52703,23 +52688,26 @@
$ws.Range("H78").Value = 90039.336
$ws.Range("J78").Value = 90000
$ws.Range("L78").Value = 270000
$ws.Range("N78").Value = -279360
# @@ -54094,25 +54082,25 @@
$ws.Range("H107").Value = 1318.7693
$ws.Range("I107").Value = 1185.25
$ws.Range("J107").Value = 1532.4
$ws.Range("K107").Value = 3555.75
$ws.Range("L107").Value = 4597.200000000001
$ws.Range("M107").Value = -1635.75
$ws.Range("N107").Value = -8437.200000000001
# @@ -54379,22 +54367,22 @@
$ws.Range("H113").Value = 361.125
$ws.Range("I113").Value = 317.16666
$ws.Range("K113").Value = 951.4999799999999
$ws.Range("M113").Value = 1218.50002
# @@ -54802,22 +54790,22 @@
$ws.Range("H122").Value = 2633.1035
$ws.Range("I122").Value = 1642.7222
$ws.Range("K122").Value = 4928.1666
$ws.Range("M122").Value = -2478.1666
# @@ -55479,25 +55467,25 @@
$ws.Range("H136").Value = 2453.3914
$ws.Range("I136").Value = 2612.0527
$ws.Range("J136").Value = 1699.75
$ws.Range("K136").Value = 7836.158100000001
$ws.Range("L136").Value = 5099.25
$ws.Range("M136").Value = -5286.158100000001
$ws.Range("N136").Value = -10199.25
